$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top; this shifts all existing rows (and their
# formatting) down by one, turning old row 1 into row 2, old row 2 into
# row 3, ... old row 31 into row 32.
$ws.Rows("1:1").Insert()

# Populate the brand-new row 1 with a simple numeric column index (0-10).
$ws.Cells.Item(1, 1).Value = 0
$ws.Cells.Item(1, 2).Value = 1
$ws.Cells.Item(1, 3).Value = 2
$ws.Cells.Item(1, 4).Value = 3
$ws.Cells.Item(1, 5).Value = 4
$ws.Cells.Item(1, 6).Value = 5
$ws.Cells.Item(1, 7).Value = 6
$ws.Cells.Item(1, 8).Value = 7
$ws.Cells.Item(1, 9).Value = 8
$ws.Cells.Item(1, 10).Value = 9
$ws.Cells.Item(1, 11).Value = 10

# Give row 1 the bold/centered/bordered "header" look (matching the style
# that used to live on the old row 1, which has now moved to row 2).
$headerRow = $ws.Range("A1:K1")
$headerRow.Font.Bold = $true
$headerRow.HorizontalAlignment = -4108
$headerRow.VerticalAlignment = -4160
$headerRow.Borders.LineStyle = 1

# Row 2 (the old header row) loses that special styling - it becomes a
# plain, unstyled row of text.
$ws.Range("A2:K2").ClearFormats()

# Row 2 keeps most of the old header captions, but the last two captions
# (thread_size / material_surface) are cleared out.
$ws.Cells.Item(2, 10).ClearContents()
$ws.Cells.Item(2, 11).ClearContents()
